$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9146023988723755
$ws.Range("B1").Value = 2.802084684371948
$ws.Range("C1").Value = 8.826982498168945
$ws.Range("D1").Value = 2.034114122390747
$ws.Range("E1").Value = 1.152721166610718
